$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert 3 new rows before row 5 (old rows 5,6,7 shift down to 8,9,10;
# formulas such as SUM(C3:C5) and the per-row $B$5 ref auto-adjust).
$ws.Rows.Item(5).Resize(3).Insert()

# Rename the existing entity rows.
$ws.Range("A3").Value = "Entity1"
$ws.Range("A4").Value = "Entity2"
$ws.Range("A8").Value = "Entity3"

# Fill in the three newly inserted attribute rows.
$ws.Range("A5").Value = "Attr1"
$ws.Range("A6").Value = "Attr2"
$ws.Range("A7").Value = "Attr3"

# Match the bold "value" column styling used by the other entity rows,
# leaving these cells without a value/formula.
$ws.Range("C5").Font.Bold = $true
$ws.Range("C6").Font.Bold = $true
$ws.Range("C7").Font.Bold = $true

# Italicize the new attribute labels.
$ws.Range("A5:A7").Font.Italic = $true

[void]$ws.Range("C9").Select()
